$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.03041816718429699, 0.5969587063387518, 0.6135877254643746, 0.7833184061825527, 0.8068171195228309, 17)
    3 = @(-0.04789084097653124, 0.5862677653873664, 0.4788442395089921, 0.6919857220412803, 0.7115722915054202, 17)
    4 = @(0.1177926702203448, 0.4865537041177135, 0.4034901151917371, 0.6352087178178029, 0.6434018492833459, 17)
    5 = @(0.0673930302034892, 0.7138065010423513, 0.7470935248604208, 0.8643457206814995, 0.8882348708881498, 17)
    6 = @(0.05303991285026294, 0.3939988076169078, 0.2884509084384082, 0.5370762594254267, 0.5508992928280517, 17)
    7 = @(0.1581559649138533, 0.592388381727795, 0.503897605990655, 0.7098574546982338, 0.7133123897002451, 17)
    8 = @(0.08165673678936797, 0.4634199137647063, 0.3661953718989556, 0.6051407868413395, 0.6180598846900753, 17)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
